$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 previously held a formula (=MODE(...)); replace it with the plain
# "Time Step" label so it matches the other header cells (B1:D1), which
# are plain shared-string text.
$ws.Range("A1").Formula = ""
$ws.Range("A1").Value = "Time Step"

# Update the Sensor 1 / Sensed State columns for rows 3-12
$ws.Range("C3").Value = 0

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("B8").Value = 1
$ws.Range("D8").Value = 1

$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 1

$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = 1
$ws.Range("D12").Value = 1
